$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Rewrite rows 2-12 in place: rows 2-3 keep content (timestamp refresh),
# a new row is inserted at 4, old rows 4-6 shift to 5-7, three new rows
# are appended at 8-9, old row 7 (now at 10) keeps content, and two more
# new rows are appended at 11-12. Hyperlinks are rebuilt from scratch so
# relationship ids stay in row order.

$ws.Hyperlinks.Delete()

$ws.Cells.Item(2, 1).Value = '2025-09-07 18:20:31'
$ws.Cells.Item(2, 2).Value = '初回 AIヘルスケア×経営支援サービス|GPT-4・LINE API活用|MVP開発パートナー募集'
$ws.Cells.Item(2, 3).Value = 'システム開発'
$ws.Cells.Item(2, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = '期限情報なし'
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5388718'
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), 'https://www.lancers.jp/work/detail/5388718')
$ws.Cells.Item(2, 6).Style = "Hyperlink"
$ws.Cells.Item(2, 7).Value = 635
$ws.Cells.Item(2, 8).Value = '🔥AI,GPT ◆開発'

$ws.Cells.Item(3, 1).Value = '2025-09-07 18:20:31'
$ws.Cells.Item(3, 2).Value = '【募集】ジャーナリングとAIをテーマにしたiOSアプリ開発'
$ws.Cells.Item(3, 3).Value = 'システム開発'
$ws.Cells.Item(3, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = '期限情報なし'
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5388502'
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), 'https://www.lancers.jp/work/detail/5388502')
$ws.Cells.Item(3, 6).Style = "Hyperlink"
$ws.Cells.Item(3, 7).Value = 378
$ws.Cells.Item(3, 8).Value = '🔥AI,Ai ◆開発 ◇アプリ'

$ws.Cells.Item(4, 1).Value = '2025-09-07 18:20:31'
$ws.Cells.Item(4, 2).Value = '【急募】Windowsサイネージシステム開発のプロフェッショナル募集'
$ws.Cells.Item(4, 3).Value = 'システム開発'
$ws.Cells.Item(4, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = '期限情報なし'
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5388877'
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), 'https://www.lancers.jp/work/detail/5388877')
$ws.Cells.Item(4, 6).Style = "Hyperlink"
$ws.Cells.Item(4, 7).Value = 125
$ws.Cells.Item(4, 8).Value = '◆開発,システム開発'

$ws.Cells.Item(5, 1).Value = '2025-09-07 18:20:31'
$ws.Cells.Item(5, 2).Value = '【急募】年間カレンダー自動作成ツールの開発依頼'
$ws.Cells.Item(5, 3).Value = 'システム開発'
$ws.Cells.Item(5, 4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = '期限情報なし'
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5388837'
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), 'https://www.lancers.jp/work/detail/5388837')
$ws.Cells.Item(5, 6).Style = "Hyperlink"
$ws.Cells.Item(5, 7).Value = 120
$ws.Cells.Item(5, 8).Value = '◆ツール,開発'

$ws.Cells.Item(6, 1).Value = '2025-09-07 18:20:31'
$ws.Cells.Item(6, 2).Value = '【急募】Instagram投稿を自動でGoogleビジネスに連携するMEOツール'
$ws.Cells.Item(6, 3).Value = 'システム開発'
$ws.Cells.Item(6, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = '期限情報なし'
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5388589'
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), 'https://www.lancers.jp/work/detail/5388589')
$ws.Cells.Item(6, 6).Style = "Hyperlink"
$ws.Cells.Item(6, 7).Value = 68
$ws.Cells.Item(6, 8).Value = '◆ツール'

$ws.Cells.Item(7, 1).Value = '2025-09-07 18:20:31'
$ws.Cells.Item(7, 2).Value = 'IB報酬を得るための高性能EA開発依頼'
$ws.Cells.Item(7, 3).Value = 'システム開発'
$ws.Cells.Item(7, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = '期限情報なし'
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5388547'
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), 'https://www.lancers.jp/work/detail/5388547')
$ws.Cells.Item(7, 6).Style = "Hyperlink"
$ws.Cells.Item(7, 7).Value = 68
$ws.Cells.Item(7, 8).Value = '◆開発'

$ws.Cells.Item(8, 1).Value = '2025-09-07 18:20:31'
$ws.Cells.Item(8, 2).Value = '【急募】Google Cloud WordPress管理画面ログイン設定'
$ws.Cells.Item(8, 3).Value = 'システム開発'
$ws.Cells.Item(8, 4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = '期限情報なし'
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5388922'
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), 'https://www.lancers.jp/work/detail/5388922')
$ws.Cells.Item(8, 6).Style = "Hyperlink"
$ws.Cells.Item(8, 7).Value = 50
$ws.Cells.Item(8, 8).Value = '◇管理 ○WordPress'

$ws.Cells.Item(9, 1).Value = '2025-09-07 18:20:31'
$ws.Cells.Item(9, 2).Value = 'LINE予約システム構築!リラクゼーションマッサージ用'
$ws.Cells.Item(9, 3).Value = 'システム開発'
$ws.Cells.Item(9, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = '期限情報なし'
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5388879'
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), 'https://www.lancers.jp/work/detail/5388879')
$ws.Cells.Item(9, 6).Style = "Hyperlink"
$ws.Cells.Item(9, 7).Value = 28
$ws.Cells.Item(9, 8).ClearContents()

$ws.Cells.Item(10, 1).Value = '2025-09-07 18:20:31'
$ws.Cells.Item(10, 2).Value = '限定公開 PR 限定公開の仕事'
$ws.Cells.Item(10, 3).Value = 'システム開発'
$ws.Cells.Item(10, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(10, 5).Value = '期限情報なし'
$ws.Cells.Item(10, 6).Value = 'https://www.lancers.jp/work/detail/5385681'
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), 'https://www.lancers.jp/work/detail/5385681')
$ws.Cells.Item(10, 6).Style = "Hyperlink"
$ws.Cells.Item(10, 7).Value = 25
$ws.Cells.Item(10, 8).ClearContents()

$ws.Cells.Item(11, 1).Value = '2025-09-07 18:20:31'
$ws.Cells.Item(11, 2).Value = '【急募】Google nonprofits設定の専門家を探しています!'
$ws.Cells.Item(11, 3).Value = 'システム開発'
$ws.Cells.Item(11, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = '期限情報なし'
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5388894'
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), 'https://www.lancers.jp/work/detail/5388894')
$ws.Cells.Item(11, 6).Style = "Hyperlink"
$ws.Cells.Item(11, 7).Value = 18
$ws.Cells.Item(11, 8).ClearContents()

$ws.Cells.Item(12, 1).Value = '2025-09-07 18:20:31'
$ws.Cells.Item(12, 2).Value = '【急募】woocommerce決済情報をスプレッドシートに反映'
$ws.Cells.Item(12, 3).Value = 'システム開発'
$ws.Cells.Item(12, 4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(12, 5).Value = '期限情報なし'
$ws.Cells.Item(12, 6).Value = 'https://www.lancers.jp/work/detail/5388904'
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), 'https://www.lancers.jp/work/detail/5388904')
$ws.Cells.Item(12, 6).Style = "Hyperlink"
$ws.Cells.Item(12, 7).Value = 10
$ws.Cells.Item(12, 8).ClearContents()

